$d = $word.ActiveDocument

$replacements = @(
    @("254×8=", "176×9="),
    @("407×5=", "481×3="),
    @("902×8=", "273×6="),
    @("164×2=", "532×5="),
    @("948×2=", "929×2="),
    @("364×8=", "783×8="),
    @("899×3=", "664×4="),
    @("337×2=", "418×5="),
    @("622×6=", "589×7="),
    @("950×2=", "815×8="),
    @("103×4=", "116×3="),
    @("318×5=", "336×8="),
    @("996×6=", "446×5="),
    @("889×9=", "509×3="),
    @("705×7=", "965×4="),
    @("549×6=", "962×7="),
    @("541×3=", "184×6="),
    @("139×2=", "164×4="),
    @("143×9=", "145×6="),
    @("435×4=", "106×3="),
    @("424×8=", "154×6="),
    @("181×9=", "199×3="),
    @("842×8=", "982×9="),
    @("452×7=", "266×8="),
    @("445×8=", "166×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Write-Output "Done"
